$d = $word.ActiveDocument

# --------------------------------------------------------------------
# Step 1: Delete the block that is being removed/relocated:
#   - old Prueba 5 tail criterio ("La aplicación debe mostrarle...")
#   - old Prueba 5 "Pasos:" + its one step (numId=8, "Continuar" sin escoger)
#   - old Prueba 6 heading paragraph ("Prueba 6")
#   - old Prueba 6 intro (its text gets merged into Prueba 5's intro instead)
#   - old Prueba 6 "Criterios de aceptación:" label (duplicate; Prueba 5 keeps its own)
# This is paragraphs 33..38 (1-indexed) in the original document.
# --------------------------------------------------------------------
$pStart = $d.Paragraphs.Item(33)
$pEnd = $d.Paragraphs.Item(38)
$delRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$delRange.Delete()

# --------------------------------------------------------------------
# Step 2: Prueba 5's intro italic text is replaced by what used to be
# Prueba 6's intro text.
# --------------------------------------------------------------------
$p31 = $d.Paragraphs.Item(31)
$r31 = $p31.Range
$r31.Find.Execute("El usuario eligió los asientos y debe seleccionar las promociones para su reserva. No desea seleccionar ninguna de ellas", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "El usuario se encuentra en el último paso de la reserva, que es el de confirmar la misma escogiendo la opción de reserva", 2)

# --------------------------------------------------------------------
# Step 3: Remove the stray trailing empty bold run after
# "Criterios de aceptación:" in (the now-merged) Prueba 5.
# --------------------------------------------------------------------
$p32 = $d.Paragraphs.Item(32)
$r32 = $p32.Range
$trailing = $d.Range($r32.End - 2, $r32.End - 1)
if ($trailing.Text -eq " ") {
    $trailing.Delete()
}

# --------------------------------------------------------------------
# Step 4: The criterio "El usuario recibe un e-mail..." (now paragraph 35)
# gains a lastRenderedPageBreak marker on its first run (it used to sit on
# the deleted "...sin escoger promociones..." step).
# --------------------------------------------------------------------
$p35 = $d.Paragraphs.Item(35)
$r35 = $p35.Range
$p35xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="Prrafodelista"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:lastRenderedPageBreak/>
    <w:t xml:space="preserve">El usuario recibe un </w:t>
  </w:r>
  <w:r>
    <w:t>e-mail notificando la confirmaci&#243;n de la reserva.</w:t>
  </w:r>
</w:p>
"@
$r35.InsertXML($p35xml)

# --------------------------------------------------------------------
# Step 5/6: Renumber the following two test headings:
#   old "Prueba 7" -> "Prueba 6"
#   old "Prueba 8" -> "Prueba 7"
# --------------------------------------------------------------------
$p40 = $d.Paragraphs.Item(40)
$p40.Range.Find.Execute("Prueba 7", $false, $false, $false, $false, $false, $true, 1, $false, "Prueba 6", 2)

$p46 = $d.Paragraphs.Item(46)
$p46.Range.Find.Execute("Prueba 8", $false, $false, $false, $false, $false, $true, 1, $false, "Prueba 7", 2)
